$d = $word.ActiveDocument

# The legend paragraph starting "Table S 2" is the document's 2nd paragraph.
$p2 = $d.Paragraphs(2).Range

# 1) Turn "...Models for bat activity..." into
#    "...Models for Pipistrellus sp. bat activity...". Scoping the Find to
#    this paragraph's own range keeps the very similarly worded Table S 1
#    legend untouched.
$null = $p2.Find.Execute("for bat activity", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "for Pipistrellus sp. bat activity", 2)

# 2) Italicise just the newly inserted species name "Pipistrellus".
$p2 = $d.Paragraphs(2).Range
$null = $p2.Find.Execute("Pipistrellus", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
$p2.Font.Italic = $true

# 3) Superscript only the "-1" exponent in "ha-1" ("ha" itself stays
#    baseline, normal text).
$p2 = $d.Paragraphs(2).Range
$null = $p2.Find.Execute("ha-1", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
$supRange = $d.Range($p2.End - 2, $p2.End)
$supRange.Font.Superscript = $true
